$d = $word.ActiveDocument

# The document's first table holds the assignment cover info:
#   Row 1: "Assignment submission" | "07 October 2022 "
#   Row 2: "Student Name"          | "  Kameshwari R"
#   Row 3: "Student Roll Number"   | "9519" + "20LCS02"
$tbl = $d.Tables.Item(1)

# 1) Submission date: 07 -> 08 October 2022
$dateCell = $tbl.Cell(1, 2)
$dateCell.Range.Find.Execute("07 October 2022", $true, $false, $false, $false, $false, $true, 1, $false, "08 October 2022", 2) | Out-Null

# 2) Student name: Kameshwari R -> Shamili N
$nameCell = $tbl.Cell(2, 2)
$nameCell.Range.Find.Execute("Kameshwari R", $true, $false, $false, $false, $false, $true, 1, $false, " Shamili N", 2) | Out-Null

# 3) Student roll number: 951920LCS02 -> 951919CS091
$rollCell = $tbl.Cell(3, 2)
$rollCell.Range.Find.Execute("20LCS02", $true, $false, $false, $false, $false, $true, 1, $false, "19CS091", 2) | Out-Null
